$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "71.614.70"
$ws.Range("E2").Value = "  +2.69%  "
$ws.Range("D3").Value = "3.658.61"
$ws.Range("E3").Value = "  +8.13%  "
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "587.59"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.91%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "179.85"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.07%  "
$ws.Range("D7").Value = "3.651.14"
$ws.Range("E7").Value = "  +8.08%  "
$ws.Range("E8").Value = "  +4.75%  "
$ws.Range("E9").Value = "  +0.08%  "
$ws.Range("E10").Value = "  +1.64%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.612"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +3.67%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "49.81"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +2.70%  "
$ws.Range("E13").Value = "  +0.26%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "682.77"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.11%  "
$ws.Range("D15").Value = "4.226.10"
$ws.Range("E15").Value = "  +7.58%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "9.01"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +4.37%  "
$ws.Range("B17").Value = "WrappedBTC"
$ws.Range("C17").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D17").Value = "71.662.69"
$ws.Range("E17").Value = "  +2.80%  "
$ws.Range("B18").Value = "WrappedEther"
$ws.Range("C18").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D18").Value = "3.638.28"
$ws.Range("E18").Value = "  +7.49%  "
$ws.Range("E19").Value = "  +1.78%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "18.25"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +3.05%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "11.65"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +3.18%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.942"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +3.11%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "6.18"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +15.22%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "17.87"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +2.96%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "103.45"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.56%  "
$ws.Range("E26").Value = "  +3.10%  "
$ws.Range("E27").Value = "  +4.96%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.20"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +4.67%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "35.35"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +5.08%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "9.20"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +5.03%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.46"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +7.65%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.26"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +10.21%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "580.68"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +4.26%  "
$ws.Range("E34").Value = "  +2.36%  "
$ws.Range("E35").Value = "  +2.38%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "60.04"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +3.64%  "
$ws.Range("D37").Value = "3.743.51"
$ws.Range("E37").Value = "  +3.71%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.999"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.06%  "
$ws.Range("E39").Value = "  +2.69%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "35.60"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.59%  "
$ws.Range("D41").Value = "0.0₃0766"
$ws.Range("E42").Value = "  +4.30%  "
$ws.Range("B43").Value = "VeChain"
$ws.Range("C43").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.0463"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +8.34%  "
$ws.Range("B44").Value = "Fetch.AI"
$ws.Range("C44").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.80"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +1.27%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.347"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +2.69%  "
$ws.Range("E46").Value = "  -0.22%  "
$ws.Range("E47").Value = "  +5.11%  "
$ws.Range("E48").Value = "  +3.72%  "
$ws.Range("E49").Value = "  +4.19%  "
$ws.Range("E50").Value = "  -0.19%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "133.82"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +2.60%  "
